# Fruta / hortaliza, semanal
# Insert a new weekly record at row 48, shifting the existing rows 48-63 down to 49-64.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 48; this shifts rows 48:63 down to 49:64
# and carries formatting (e.g. the date style on column D) along.
$ws.Rows("48:48").Insert()

# Populate the newly inserted row 48 with the new record.
$ws.Range("A48").Value = 10
$ws.Range("B48").Value = "Vega Modelo de Temuco"
$ws.Range("C48").Value = "La Araucanía"
$ws.Range("D48").Value = 44837
$ws.Range("E48").Value = 9
$ws.Range("F48").Value = 100112010
$ws.Range("G48").Value = "Achicoria"
$ws.Range("H48").Value = "Sin especificar"
$ws.Range("I48").Value = "Primera"
$ws.Range("J48").Value = 300
$ws.Range("K48").Value = 8000
$ws.Range("L48").Value = 8000
$ws.Range("M48").Value = 8000
$ws.Range("N48").Value = "$/caja 18 unidades"
$ws.Range("O48").Value = "Región del Maule"
$ws.Range("P48").Value = 444
$ws.Range("Q48").Value = 18
$ws.Range("R48").Value = "Hortaliza"
